{"js": "// Simplify the tensile report RESULTS table:\n// - Rename \"Elongation at Fracture (A)\" row to \"Elongation A5 (L1-L0)\" and\n//   swap its placeholders {{A}} / {{A_uncertainty}} for {{A5_value}} / {{A5_uncertainty}}.\n// - Remove the \"Uniform Elongation (Ag)\" row entirely.\n// - Remove the \"Stress Rate at Yield\" / \"Strain Rate at Yield\" /\n//   \"Stress Rate at Rm\" / \"Strain Rate at Rm\" rows entirely.\n\nconst ROWS_TO_REMOVE = new Set([\n  \"Uniform Elongation (Ag)\",\n  \"Stress Rate at Yield\",\n  \"Strain Rate at Yield\",\n  \"Stress Rate at Rm\",\n  \"Strain Rate at Rm\",\n]);\n\n// Locate the RESULTS table: the one whose header row starts with \"Parameter\".\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet resultsTable = null;\nfor (let t = 0; t < tables.items.length; t++) {\n  const tbl = tables.items[t];\n  tbl.rows.load(\"items\");\n  await context.sync();\n  const firstRow = tbl.rows.items[0];\n  firstRow.cells.load(\"items\");\n  await context.sync();\n  const firstCell = firstRow.cells.items[0];\n  firstCell.load(\"value\");\n  await context.sync();\n  if ((firstCell.value || \"\").trim() === \"Parameter\") {\n    resultsTable = tbl;\n    break;\n  }\n}\n\nif (!resultsTable) {\n  throw new Error(\"Could not locate the RESULTS table\");\n}\n\n// Rename the \"Elongation at Fracture (A)\" row + swap its placeholders.\nresultsTable.rows.load(\"items\");\nawait context.sync();\n\nfor (const row of resultsTable.rows.items) {\n  row.cells.load(\"items\");\n  await context.sync();\n  const firstCell = row.cells.items[0];\n  firstCell.load(\"value\");\n  await context.sync();\n  if ((firstCell.value || \"\").trim() === \"Elongation at Fracture (A)\") {\n    const valueCell = row.cells.items[1];\n    const uncertaintyCell = row.cells.items[2];\n    valueCell.load(\"value\");\n    uncertaintyCell.load(\"value\");\n    await context.sync();\n\n    firstCell.value = \"Elongation A5 (L1-L0)\";\n    valueCell.value = valueCell.value.replace(\"{{A}}\", \"{{A5_value}}\");\n    uncertaintyCell.value = uncertaintyCell.value.replace(\"{{A_uncertainty}}\", \"{{A5_uncertainty}}\");\n    await context.sync();\n    break;\n  }\n}\n\n// Remove the obsolete rows. Re-query fresh each time (instead of reusing a\n// stale `items` array) so row anchors stay valid across the deletes.\nasync function findRowToRemove() {\n  resultsTable.rows.load(\"items\");\n  await context.sync();\n  for (const row of resultsTable.rows.items) {\n    row.cells.load(\"items\");\n    await context.sync();\n    const firstCell = row.cells.items[0];\n    firstCell.load(\"value\");\n    await context.sync();\n    if (ROWS_TO_REMOVE.has((firstCell.value || \"\").trim())) {\n      return row;\n    }\n  }\n  return null;\n}\n\nlet rowToRemove;\nwhile ((rowToRemove = await findRowToRemove())) {\n  rowToRemove.delete();\n  await context.sync();\n}\n", "ps1": "# Simplify the tensile report RESULTS table:\n# - Rename \"Elongation at Fracture (A)\" row to \"Elongation A5 (L1-L0)\" and\n#   swap its placeholders {{A}} / {{A_uncertainty}} for {{A5_value}} / {{A5_uncertainty}}.\n# - Remove the \"Uniform Elongation (Ag)\" row entirely.\n# - Remove the \"Stress Rate at Yield\" / \"Strain Rate at Yield\" /\n#   \"Stress Rate at Rm\" / \"Strain Rate at Rm\" rows entirely.\n\n$d = $word.ActiveDocument\n\n# Locate the RESULTS table: the one whose header row starts with \"Parameter\".\n$resultsTable = $null\nfor ($t = 1; $t -le $d.Tables.Count; $t++) {\n    $candidate = $d.Tables.Item($t)\n    $header = $candidate.Cell(1, 1).Range.Text.TrimEnd([char]13, [char]7)\n    if ($header -eq \"Parameter\") {\n        $resultsTable = $candidate\n        break\n    }\n}\n\n$rowsToRemove = @(\n    \"Uniform Elongation (Ag)\",\n    \"Stress Rate at Yield\",\n    \"Strain Rate at Yield\",\n    \"Stress Rate at Rm\",\n    \"Strain Rate at Rm\"\n)\n\n# Walk bottom-to-top so deleting a row never shifts the index of a row we\n# still need to inspect/delete.\nfor ($r = $resultsTable.Rows.Count; $r -ge 1; $r--) {\n    $label = $resultsTable.Cell($r, 1).Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($label -eq \"Elongation at Fracture (A)\") {\n        $resultsTable.Cell($r, 1).Range.Text = \"Elongation A5 (L1-L0)\"\n\n        $valueText = $resultsTable.Cell($r, 2).Range.Text.TrimEnd([char]13, [char]7)\n        $resultsTable.Cell($r, 2).Range.Text = $valueText.Replace(\"{{A}}\", \"{{A5_value}}\")\n\n        $uncertaintyText = $resultsTable.Cell($r, 3).Range.Text.TrimEnd([char]13, [char]7)\n        $resultsTable.Cell($r, 3).Range.Text = $uncertaintyText.Replace(\"{{A_uncertainty}}\", \"{{A5_uncertainty}}\")\n    }\n    elseif ($rowsToRemove -contains $label) {\n        $resultsTable.Rows.Item($r).Delete()\n    }\n}\n"}
